# Insert a new data row at row 362 (this pushes the former rows 362-388 down
# to become rows 363-389), then populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("362:362").Insert()

$ws.Range("A362").Value2 = 11
$ws.Range("B362").Value2 = "Vega Monumental Concepción"
$ws.Range("C362").Value2 = "Bíobío"
$ws.Range("D362").Value2 = 45013
$ws.Range("E362").Value2 = 8
$ws.Range("F362").Value2 = 100112009
$ws.Range("G362").Value2 = "Acelga"
$ws.Range("H362").Value2 = "Sin especificar"
$ws.Range("I362").Value2 = "Primera"
$ws.Range("J362").Value2 = 270
$ws.Range("K362").Value2 = 600
$ws.Range("L362").Value2 = 650
$ws.Range("M362").Value2 = 622
$ws.Range("N362").Value2 = "`$/atado 0,5 a 1 kilo"
$ws.Range("O362").Value2 = "Región de Ñuble"
$ws.Range("P362").Value2 = 622
$ws.Range("Q362").Value2 = 1
$ws.Range("R362").Value2 = "Hortaliza"
